$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 29: shareclass_dist_status
$ws.Range("A29").Value = "shareclass_dist_status"
$ws.Range("B29").Value = "shareclass_dist_status"
$ws.Range("C29").Value = "https://am.jpmorgan.com/"
$ws.Hyperlinks.Add($ws.Range("C29"), "https://am.jpmorgan.com/")
$ws.Range("C29").Style = $ws.Range("C28").Style

# New row 30: shareclass_assets
$ws.Range("A30").Value = "shareclass_assets"
$ws.Range("B30").Value = "shareclass_assets"
$ws.Range("C30").Value = "https://am.jpmorgan.com/"
$ws.Hyperlinks.Add($ws.Range("C30"), "https://am.jpmorgan.com/")
$ws.Range("C30").Style = $ws.Range("C28").Style

# Update selection to match the post-edit state
$ws.Range("E31").Select()
